$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain Text so numeric-looking strings
# (e.g. "1.000", "30.90", "0.000006940") keep their exact literal formatting
# instead of being auto-coerced to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "26.499.92"
$ws.Cells.Item(2, 5).Value = "  +0.61%  "
$ws.Cells.Item(3, 4).Value = "1.727.91"
$ws.Cells.Item(3, 5).Value = "  +0.31%  "
$ws.Cells.Item(4, 4).Value = "0.9997"
$ws.Cells.Item(4, 5).Value = "  -0.04%  "
$ws.Cells.Item(5, 4).Value = "245.39"
$ws.Cells.Item(5, 5).Value = "  +2.53%  "
$ws.Cells.Item(6, 4).Value = "1.000"
$ws.Cells.Item(6, 5).Value = "  -0.07%  "
$ws.Cells.Item(7, 4).Value = "0.4807"
$ws.Cells.Item(7, 5).Value = "  +1.94%  "
$ws.Cells.Item(8, 4).Value = "0.2669"
$ws.Cells.Item(8, 5).Value = "  +1.26%  "
$ws.Cells.Item(9, 4).Value = "0.06223"
$ws.Cells.Item(9, 5).Value = "  +0.14%  "
$ws.Cells.Item(10, 4).Value = "1.725.25"
$ws.Cells.Item(10, 5).Value = "  +0.12%  "
$ws.Cells.Item(11, 4).Value = "0.07155"
$ws.Cells.Item(11, 5).Value = "  +1.14%  "
$ws.Cells.Item(12, 4).Value = "15.65"
$ws.Cells.Item(12, 5).Value = "  +2.10%  "
$ws.Cells.Item(13, 4).Value = "0.6161"
$ws.Cells.Item(13, 5).Value = "  +4.18%  "
$ws.Cells.Item(14, 5).Value = "  +2.80%  "
$ws.Cells.Item(15, 4).Value = "77.17"
$ws.Cells.Item(15, 5).Value = "  +1.18%  "
$ws.Cells.Item(16, 5).Value = "  -0.06%  "
$ws.Cells.Item(17, 4).Value = "26.507.56"
$ws.Cells.Item(18, 4).Value = "1.000"
$ws.Cells.Item(18, 5).Value = "  -0.11%  "
$ws.Cells.Item(19, 4).Value = "0.000006940"
$ws.Cells.Item(19, 5).Value = "  +2.33%  "
$ws.Cells.Item(20, 4).Value = "11.64"
$ws.Cells.Item(20, 5).Value = "  +0.65%  "
$ws.Cells.Item(21, 4).Value = "1.946.56"
$ws.Cells.Item(21, 5).Value = "  +0.34%  "
$ws.Cells.Item(22, 4).Value = "4.527"
$ws.Cells.Item(22, 5).Value = "  -0.61%  "
$ws.Cells.Item(23, 4).Value = "8.960"
$ws.Cells.Item(23, 5).Value = "  +2.40%  "
$ws.Cells.Item(24, 4).Value = "5.281"
$ws.Cells.Item(24, 5).Value = "  -1.08%  "
$ws.Cells.Item(25, 4).Value = "136.68"
$ws.Cells.Item(25, 5).Value = "  +1.45%  "
$ws.Cells.Item(26, 5).Value = "  +0.80%  "
$ws.Cells.Item(27, 5).Value = "  +1.98%  "
$ws.Cells.Item(28, 4).Value = "1.405"
$ws.Cells.Item(28, 5).Value = "  -0.27%  "
$ws.Cells.Item(29, 4).Value = "106.86"
$ws.Cells.Item(29, 5).Value = "  -1.57%  "
$ws.Cells.Item(30, 4).Value = "3.975"
$ws.Cells.Item(30, 5).Value = "  -0.83%  "
$ws.Cells.Item(31, 4).Value = "0.08032"
$ws.Cells.Item(31, 5).Value = "  +3.79%  "
$ws.Cells.Item(32, 5).Value = "  +0.62%  "
$ws.Cells.Item(33, 4).Value = "0.04562"
$ws.Cells.Item(33, 5).Value = "  +2.24%  "
$ws.Cells.Item(34, 2).Value = "Frax"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(34, 4).Value = "0.9996"
$ws.Cells.Item(34, 5).Value = "  -0.08%  "
$ws.Cells.Item(35, 2).Value = "HuobiToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(35, 4).Value = "2.616"
$ws.Cells.Item(35, 5).Value = "  +0.04%  "
$ws.Cells.Item(36, 2).Value = "ImmutableX"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(36, 4).Value = "0.6364"
$ws.Cells.Item(36, 5).Value = "  +2.66%  "
$ws.Cells.Item(37, 2).Value = "ARBITRUM"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(37, 4).Value = "0.9929"
$ws.Cells.Item(37, 5).Value = "  +1.76%  "
$ws.Cells.Item(38, 2).Value = "TrustWalletToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(38, 4).Value = "0.9339"
$ws.Cells.Item(38, 5).Value = "  +1.13%  "
$ws.Cells.Item(39, 2).Value = "RenderToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(39, 4).Value = "2.095"
$ws.Cells.Item(39, 5).Value = "  +10.17%  "
$ws.Cells.Item(40, 2).Value = "MXToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(40, 4).Value = "2.426"
$ws.Cells.Item(40, 5).Value = "  +0.08%  "
$ws.Cells.Item(41, 2).Value = "Quant"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(41, 4).Value = "105.45"
$ws.Cells.Item(41, 5).Value = "  -8.33%  "
$ws.Cells.Item(42, 2).Value = "PaxDollar"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(42, 4).Value = "1.002"
$ws.Cells.Item(42, 5).Value = "  +0.10%  "
$ws.Cells.Item(43, 2).Value = "VeChain"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(43, 4).Value = "0.01503"
$ws.Cells.Item(43, 5).Value = "  +1.57%  "
$ws.Cells.Item(44, 2).Value = "FraxShare"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(44, 4).Value = "5.588"
$ws.Cells.Item(44, 5).Value = "  +3.47%  "
$ws.Cells.Item(45, 2).Value = "TheSandbox"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(45, 4).Value = "0.3907"
$ws.Cells.Item(45, 5).Value = "  +2.37%  "
$ws.Cells.Item(46, 2).Value = "Aptos"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(46, 4).Value = "6.925"
$ws.Cells.Item(46, 5).Value = "  +10.78%  "
$ws.Cells.Item(47, 2).Value = "Algorand"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(47, 4).Value = "0.1185"
$ws.Cells.Item(47, 5).Value = "  +1.87%  "
$ws.Cells.Item(48, 2).Value = "Cronos"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(48, 4).Value = "0.05331"
$ws.Cells.Item(48, 5).Value = "  +0.81%  "
$ws.Cells.Item(49, 2).Value = "Elrond"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(49, 4).Value = "30.90"
$ws.Cells.Item(49, 5).Value = "  +0.91%  "
$ws.Cells.Item(50, 2).Value = "EnergySwap"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(50, 4).Value = "7.886"
$ws.Cells.Item(50, 5).Value = "  +2.61%  "
$ws.Cells.Item(51, 2).Value = "NEARProtocol"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(51, 4).Value = "1.269"
$ws.Cells.Item(51, 5).Value = "  +4.21%  "
